$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview sheet (also mirrored by de-de!H2 via the
# same underlying shared string), updated to the new report-generation timestamp.
$wsOverview.Range("G2").Value = "2016-08-25 17:07:43"
$wsDeDe.Range("H2").Value = "2016-08-25 17:07:43"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-08-25 17:07:38"
$wsZhCn.Range("K2").Value = "2016-08-25 17:08:21"

# de-de sheet: Correspond Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-08-25 17:08:28"
